$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "full_teste" token row (row 9) below the existing token rows,
# mirroring the layout of the previous rows (token name in D, an "x" mark
# in whichever permission columns apply, and the token value in J).
# Copy the formatting from the row right above it (row 8) first so the new
# row keeps the same look (borders/alignment/fonts) as rows 6-8.
$ws.Range("D8:J8").Copy()
$ws.Range("D9:J9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D9").Value = "full_teste"
$ws.Range("E9").Value = "x"
$ws.Range("F9").Value = "x"
$ws.Range("G9").Value = "x"
$ws.Range("H9").Value = "x"
$ws.Range("I9").Value = "x"
$ws.Range("J9").Value = "5a7J4zG40xIUAZe"

# Move the active selection the way it ended up after entering the new row.
[void]$ws.Range("D12").Select()
